$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rangeB1 = $ws.Range("B$r1")
    $rangeB2 = $ws.Range("B$r2")
    $rangeR1 = $ws.Range("F${r1}:AC${r1}")
    $rangeR2 = $ws.Range("F${r2}:AC${r2}")

    $b1 = $rangeB1.Value2
    $b2 = $rangeB2.Value2
    $v1 = $rangeR1.Value2
    $v2 = $rangeR2.Value2

    $rangeB1.Value2 = $b2
    $rangeB2.Value2 = $b1
    $rangeR1.Value2 = $v2
    $rangeR2.Value2 = $v1
}

Swap-Rows 25 26
Swap-Rows 84 85
Swap-Rows 226 227
Swap-Rows 229 230
Swap-Rows 263 264
Swap-Rows 276 278
Swap-Rows 285 286
Swap-Rows 313 314
Swap-Rows 328 329
Swap-Rows 352 353
Swap-Rows 387 388
Swap-Rows 411 412
Swap-Rows 437 438
Swap-Rows 500 501
